$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.674.25"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.669.39"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'599.20"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "156.47"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.617"
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("E9").Value = "  +4.49%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "29.27"
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").Value = "3.149.39"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "65.517.27"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.669.74"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "351.89"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "69.63"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  +4.89%  "
$ws.Range("D25").Value = "9.61"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value = "  -5.49%  "
$ws.Range("D29").Value = "'8.00"
$ws.Range("E29").Value = "  -4.88%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("D32").Value = "527.72"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").Value = "20.54"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "158.59"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "163.52"
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("E44").Value = "  +1.68%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "22.72"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").Value = "0.637"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("E49").Value = "  +14.28%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "20.09"
$ws.Range("E51").Value = "  -4.02%  "
